$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "HK_G_acc_G"

$values = @(
50.070323488045013,
49.929676511954995,
50.210970464135016,
49.507735583684955,
49.507735583684955,
49.367088607594937,
51.61744022503516,
51.336146272855132,
51.61744022503516,
51.898734177215189,
50.351617440225041,
50.914205344585092,
51.758087201125178,
51.61744022503516,
52.180028129395218,
50.492264416315045,
50.773558368495074,
49.929676511954995,
50.632911392405063,
51.054852320675103,
50.914205344585092,
48.804500703234879,
48.38255977496484,
48.52320675105485,
50.632911392405063,
50.914205344585092,
51.195499296765121,
52.180028129395218,
51.47679324894515,
51.758087201125178,
48.945147679324897,
49.367088607594937,
49.789029535864984,
50.632911392405063,
50.070323488045013,
53.586497890295362,
48.945147679324897,
49.929676511954995,
49.929676511954995,
51.47679324894515,
51.61744022503516,
51.758087201125178,
51.61744022503516,
50.351617440225041,
50.351617440225041,
49.226441631504926,
51.898734177215189,
50.632911392405063
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
